# Refresh the cryptos price table (generated from the upstream GitHub Actions diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.609.79"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.664.35"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'0.9989"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").Value = "'302.45"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.3833"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").Value = "'51.26"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "'0.3604"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "'1.249"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("D11").Value = "'0.08194"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "'0.9987"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "'22.51"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").Value = "'6.531"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("D15").Value = "'7.534"
$ws.Range("E15").Value = "  +4.51%  "
$ws.Range("D16").Value = "'0.00001226"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "1.659.71"
$ws.Range("E17").Value = "  +3.39%  "
$ws.Range("D18").Value = "'97.65"
$ws.Range("E18").Value = "  +3.80%  "
$ws.Range("D19").Value = "'0.06982"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("D20").Value = "'6.842"
$ws.Range("E20").Value = "  +4.94%  "
$ws.Range("D21").Value = "'17.77"
$ws.Range("E21").Value = "  +3.72%  "
$ws.Range("D22").Value = "'0.9987"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").Value = "'12.77"
$ws.Range("E23").Value = "  +3.86%  "
$ws.Range("D24").Value = "23.629.29"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "'2.523"
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'3.022"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "'21.28"
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("D28").Value = "'152.52"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'5.240"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("D30").Value = "'134.15"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'7.260"
$ws.Range("E31").Value = "  +11.68%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.846.41"
$ws.Range("E32").Value = "  +3.37%  "
$ws.Range("E33").Value = "  +7.12%  "
$ws.Range("D34").Value = "'11.96"
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'0.02808"
$ws.Range("E36").Value = "  +4.00%  "
$ws.Range("D37").Value = "'6.179"
$ws.Range("E37").Value = "  +6.41%  "
$ws.Range("D38").Value = "'0.2507"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'0.07068"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'13.44"
$ws.Range("E41").Value = "  +12.74%  "
$ws.Range("D42").Value = "'0.7057"
$ws.Range("D43").Value = "'1.336"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "'16.05"
$ws.Range("E44").Value = "  +5.50%  "
$ws.Range("D45").Value = "'0.6573"
$ws.Range("E45").Value = "  +4.58%  "
$ws.Range("D46").Value = "'2.320"
$ws.Range("E46").Value = "  +3.77%  "
$ws.Range("D48").Value = "'3.967"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'0.07968"
$ws.Range("E49").Value = "  +1.56%  "
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "'1.195"
$ws.Range("E51").Value = "  +2.91%  "
